$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters corresponding to each value in a row (A and G/H are unchanged)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# New values for rows 2..25, one array per row, values ordered per $cols
$data = @(
    @(1.02,1.03554996892312,1.037993198678242,1.044620814953643,1.055707090947374,1.034684832415623,1.040663132702154,1.040782499662942,1.04739135119426,1.058446809463551,1.017422043476474),
    @(1.02,1.036744317178569,1.038882382051001,1.045701760126028,1.056919383875835,1.03491832844936,1.041499840663056,1.04148139530513,1.048282868076084,1.059471580434422,1.017706497846563),
    @(1.02,1.037516951268669,1.039457338508437,1.046401376487236,1.057704130960776,1.035067918078169,1.042040553845411,1.041932591279214,1.048859318657623,1.060134419458266,1.017890148103943),
    @(1.02,1.037841723104329,1.039698953752248,1.046695537421018,1.058034114626011,1.035130447101403,1.042267704655212,1.042122026515821,1.049101558671348,1.060413016618033,1.017967256519623),
    @(1.02,1.037896251222055,1.039739516371103,1.046744930838021,1.058089524825794,1.035140924983759,1.042305834623979,1.042153819024377,1.04914222600721,1.060459790796146,1.0179801976114),
    @(1.02,1.037521291054726,1.039460567361618,1.04640530691302,1.057708539918595,1.035068755002344,1.042043589692976,1.041935123494172,1.048862555871676,1.060138142323102,1.017891178815816),
    @(1.02,1.03595364470517,1.038293786401309,1.044986090253185,1.056116727023798,1.034764053736911,1.040946046039756,1.041018909596197,1.047692731067122,1.058793189284555,1.017518261318024),
    @(1.02,1.033189725398747,1.036234650659496,1.042486529684803,1.053314089485722,1.034215651999504,1.039006689251297,1.039396463166534,1.045628091752955,1.056421202005274,1.016857978970364),
    @(1.02,1.031345975155068,1.034859766450972,1.040820956806505,1.051447160725998,1.033842321908519,1.037710135116901,1.038309443898242,1.044249418639351,1.05483846132683,1.016415657495444),
    @(1.02,1.030547315669603,1.034263914242605,1.040099921080591,1.050639096828466,1.033678829784512,1.037147833979185,1.037837465711272,1.043651891795831,1.05415276590344,1.016223618299128),
    @(1.02,1.030250610089297,1.034042509568874,1.039832120270889,1.05033899362708,1.033617825119346,1.036938836155783,1.037661957081278,1.043429859511845,1.053898012837592,1.016152209295498),
    @(1.02,1.030314256647359,1.034090005193689,1.039889563390601,1.050403364650193,1.033630923340315,1.036983672988589,1.037699613112393,1.043477490020988,1.053952660758405,1.016167530262948),
    @(1.02,1.030522790883948,1.034245614471766,1.040077784100963,1.050614289240558,1.033673792760368,1.037130560897579,1.03782296209438,1.043633540262974,1.054131709071015,1.016217717183149),
    @(1.02,1.03065126929538,1.03434148004099,1.040193756228294,1.050744253114414,1.033700169388446,1.037221045573714,1.037898935619154,1.043729676771809,1.054242019281096,1.016248628775789),
    @(1.02,1.031398972965371,1.034899300214767,1.040868813006283,1.051500796047079,1.033853133596817,1.037747434447649,1.038340740248973,1.044289062786713,1.054883961003621,1.016428391710035),
    @(1.02,1.031867904858167,1.035249066774679,1.041292302418019,1.051975442242019,1.033948591770279,1.038077386565144,1.038617526199534,1.044639802121325,1.055286537159879,1.016541015199518),
    @(1.02,1.032141395774199,1.035453029739068,1.041539332958601,1.052252327153258,1.034004093615126,1.03826975678266,1.038778846193321,1.044844329238103,1.055521318551354,1.016606657285389),
    @(1.02,1.03223464423799,1.035522567388946,1.041623566798757,1.052346743198705,1.03402298824724,1.038335335629171,1.038833831038861,1.044914058712166,1.055601367134661,1.016629031189109),
    @(1.02,1.031817595934127,1.035211545306152,1.041246864333306,1.051924513919484,1.033938368352498,1.038041994655265,1.038587842578136,1.044602176622918,1.055243348124903,1.016528936865806),
    @(1.02,1.030461384089225,1.034199793594696,1.040022357142212,1.050552175902164,1.033661176417463,1.037087309782465,1.037786644281138,1.043587589690687,1.054078985342425,1.016202940516756),
    @(1.02,1.029608401674994,1.033563209804986,1.03925259885356,1.049689607730646,1.033485295462683,1.036486284777942,1.037281770957603,1.042949190793606,1.05334658384837,1.015997527466184),
    @(1.02,1.030060610929644,1.033900718335064,1.039660649651514,1.050146845906336,1.033578684992448,1.036804973544343,1.037549521082655,1.043287664781412,1.053734874586694,1.016106463210854),
    @(1.02,1.031840328461404,1.035228499797176,1.041267395797585,1.051947526138058,1.033942988422791,1.038057986999465,1.038601255713922,1.044619178130557,1.055262863497362,1.016534394697248),
    @(1.02,1.033904457092182,1.036767359271846,1.043132580431483,1.054038369074772,1.03435878846002,1.039508698603677,1.039816851520411,1.04616224241508,1.057034662298065,1.017029052849228)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $row).Value = $rowValues[$j]
    }
}
